$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 46851.668
$ws.Range("J105").Value = 46851.668
$ws.Range("L105").Value = 46851.668
$ws.Range("N105").Value = -53839.668

$ws.Range("H113").Value = 1499.5
$ws.Range("I113").Value = 1499
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1499
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1755
$ws.Range("N113").Value = -8008

$ws.Range("H141").Value = 3360.3333
$ws.Range("I141").Value = 1574.1923
$ws.Range("J141").Value = 14970.25
$ws.Range("K141").Value = 4722.5769
$ws.Range("L141").Value = 44910.75
$ws.Range("M141").Value = 457.4231
$ws.Range("N141").Value = -55270.75


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 7249.75
$ws.Range("J10").Value = 9499.5
$ws.Range("L10").Value = 9499.5
$ws.Range("N10").Value = -9839.5

$ws.Range("H24").Value = 29500
$ws.Range("J24").Value = 29500
$ws.Range("L24").Value = 29500
$ws.Range("N24").Value = -30248

$ws.Range("H32").Value = 3149.3635
$ws.Range("I32").Value = 2371.0334
$ws.Range("K32").Value = 2371.0334
$ws.Range("M32").Value = -2084.0334

$ws.Range("H97").Value = 566.6
$ws.Range("I97").Value = 542.2222
$ws.Range("K97").Value = 542.2222
$ws.Range("M97").Value = -46.22220000000004

$ws.Range("H100").Value = 29500
$ws.Range("J100").Value = 29500
$ws.Range("L100").Value = 29500
$ws.Range("N100").Value = -31664

$ws.Range("H122").Value = 5037499.5
$ws.Range("I122").Value = 5037499.5
$ws.Range("K122").Value = 15112498.5
$ws.Range("M122").Value = -15110048.5


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 414.6
$ws.Range("I5").Value = 492.42856
$ws.Range("J5").Value = 233
$ws.Range("K5").Value = 492.42856
$ws.Range("L5").Value = 233
$ws.Range("M5").Value = -379.42856
$ws.Range("N5").Value = -459

$ws.Range("H81").Value = 37500
$ws.Range("J81").Value = 35000
$ws.Range("L81").Value = 35000
$ws.Range("N81").Value = -37122

$ws.Range("H84").Value = 37500
$ws.Range("J84").Value = 35000
$ws.Range("L84").Value = 105000
$ws.Range("N84").Value = -115608

$ws.Range("H100").Value = 25410.5
$ws.Range("J100").Value = 25410.5
$ws.Range("L100").Value = 25410.5
$ws.Range("N100").Value = -27574.5

$ws.Range("H107").Value = 1790.7858
$ws.Range("I107").Value = 1736.2307
$ws.Range("K107").Value = 1736.2307
$ws.Range("M107").Value = 183.7692999999999

$ws.Range("H134").Value = 1931.125
$ws.Range("I134").Value = 1423.1333
$ws.Range("J134").Value = 2777.7778
$ws.Range("K134").Value = 4269.3999
$ws.Range("L134").Value = 8333.3334
$ws.Range("M134").Value = -1734.3999
$ws.Range("N134").Value = -13403.3334


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2658.0571
$ws.Range("I31").Value = 1715.1072
$ws.Range("J31").Value = 6429.857
$ws.Range("K31").Value = 1715.1072
$ws.Range("L31").Value = 6429.857
$ws.Range("M31").Value = -1420.1072
$ws.Range("N31").Value = -7019.857

$ws.Range("H34").Value = 2658.0571
$ws.Range("I34").Value = 1715.1072
$ws.Range("J34").Value = 6429.857
$ws.Range("K34").Value = 1715.1072
$ws.Range("L34").Value = 6429.857
$ws.Range("M34").Value = -1513.1072
$ws.Range("N34").Value = -6833.857

$ws.Range("H99").Value = 9791.156000000001
$ws.Range("I99").Value = 5576.9473
$ws.Range("K99").Value = 5576.9473
$ws.Range("M99").Value = -4078.9473

$ws.Range("H106").Value = 33999
$ws.Range("J106").Value = 33999
$ws.Range("L106").Value = 33999
$ws.Range("N106").Value = -36523

$ws.Range("H122").Value = 1319.8
$ws.Range("I122").Value = 1266.3334
$ws.Range("K122").Value = 3799.0002
$ws.Range("M122").Value = -1349.0002

$ws.Range("H126").Value = 9791.156000000001
$ws.Range("I126").Value = 5576.9473
$ws.Range("K126").Value = 16730.8419
$ws.Range("M126").Value = -14260.8419

$ws.Range("H132").Value = 2145.9443
$ws.Range("I132").Value = 1795.3572
$ws.Range("K132").Value = 5386.071599999999
$ws.Range("M132").Value = -2856.071599999999


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 423.2
$ws.Range("I122").Value = 232
$ws.Range("J122").Value = 710
$ws.Range("K122").Value = 2088
$ws.Range("L122").Value = 6390
$ws.Range("M122").Value = 362
$ws.Range("N122").Value = -11290


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 4500
$ws.Range("J12").Value = 5000
$ws.Range("L12").Value = 5000
$ws.Range("N12").Value = -5280

$ws.Range("H18").Value = 2012720
$ws.Range("J18").Value = 15900
$ws.Range("L18").Value = 15900
$ws.Range("N18").Value = -16486

$ws.Range("H20").Value = 30036.8
$ws.Range("J20").Value = 30036.8
$ws.Range("L20").Value = 30036.8
$ws.Range("N20").Value = -30526.8

$ws.Range("H55").Value = 8000
$ws.Range("J55").Value = 8000
$ws.Range("L55").Value = 8000
$ws.Range("N55").Value = -8654


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1915.4445
$ws.Range("J16").Value = 1898.6666
$ws.Range("L16").Value = 1898.6666
$ws.Range("N16").Value = -2238.6666

$ws.Range("H40").Value = 3602.2
$ws.Range("I40").Value = 2707.8333
$ws.Range("J40").Value = 4943.75
$ws.Range("K40").Value = 2707.8333
$ws.Range("L40").Value = 4943.75
$ws.Range("M40").Value = -2571.8333
$ws.Range("N40").Value = -5215.75

$ws.Range("H41").Value = 26166.666
$ws.Range("J41").Value = 26166.666
$ws.Range("L41").Value = 26166.666
$ws.Range("N41").Value = -27042.666

$ws.Range("H47").Value = 23666.666
$ws.Range("J47").Value = 23666.666
$ws.Range("L47").Value = 23666.666
$ws.Range("N47").Value = -24646.666

$ws.Range("H52").Value = 23666.666
$ws.Range("J52").Value = 23666.666
$ws.Range("L52").Value = 23666.666
$ws.Range("N52").Value = -24132.666

$ws.Range("H122").Value = 2333
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 753748
$ws.Range("J5").Value = 4997.3335
$ws.Range("L5").Value = 4997.3335
$ws.Range("N5").Value = -5221.3335

$ws.Range("H11").Value = 19994.5
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 19994.5
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 19994.5
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = -20278.5

$ws.Range("H12").Value = 5333.3335
$ws.Range("J12").Value = 6000
$ws.Range("L12").Value = 6000
$ws.Range("N12").Value = -6284

$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").Value = ""

$ws.Range("H103").Value = 14767.333
$ws.Range("J103").Value = 14767.333
$ws.Range("L103").Value = 14767.333
$ws.Range("N103").Value = -17111.333

